$d = $word.ActiveDocument

$d.Content.Find.Execute("540÷6=90, 0", $true, $false, $false, $false, $false, $true, 1, $false, "126÷4=31, 2", 2) | Out-Null
$d.Content.Find.Execute("711÷3=237, 0", $true, $false, $false, $false, $false, $true, 1, $false, "707÷7=101, 0", 2) | Out-Null
$d.Content.Find.Execute("918÷9=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "804÷5=160, 4", 2) | Out-Null
$d.Content.Find.Execute("549÷9=61, 0", $true, $false, $false, $false, $false, $true, 1, $false, "958÷6=159, 4", 2) | Out-Null
$d.Content.Find.Execute("407÷3=135, 2", $true, $false, $false, $false, $false, $true, 1, $false, "913÷2=456, 1", 2) | Out-Null
$d.Content.Find.Execute("180÷7=25, 5", $true, $false, $false, $false, $false, $true, 1, $false, "805÷9=89, 4", 2) | Out-Null
$d.Content.Find.Execute("891÷3=297, 0", $true, $false, $false, $false, $false, $true, 1, $false, "995÷9=110, 5", 2) | Out-Null
$d.Content.Find.Execute("498÷6=83, 0", $true, $false, $false, $false, $false, $true, 1, $false, "326÷9=36, 2", 2) | Out-Null
$d.Content.Find.Execute("687÷2=343, 1", $true, $false, $false, $false, $false, $true, 1, $false, "116÷5=23, 1", 2) | Out-Null
$d.Content.Find.Execute("845÷7=120, 5", $true, $false, $false, $false, $false, $true, 1, $false, "895÷3=298, 1", 2) | Out-Null
$d.Content.Find.Execute("591÷6=98, 3", $true, $false, $false, $false, $false, $true, 1, $false, "256÷7=36, 4", 2) | Out-Null
$d.Content.Find.Execute("603÷4=150, 3", $true, $false, $false, $false, $false, $true, 1, $false, "677÷9=75, 2", 2) | Out-Null
$d.Content.Find.Execute("222÷5=44, 2", $true, $false, $false, $false, $false, $true, 1, $false, "815÷9=90, 5", 2) | Out-Null
$d.Content.Find.Execute("263÷4=65, 3", $true, $false, $false, $false, $false, $true, 1, $false, "274÷3=91, 1", 2) | Out-Null
$d.Content.Find.Execute("571÷9=63, 4", $true, $false, $false, $false, $false, $true, 1, $false, "603÷8=75, 3", 2) | Out-Null
$d.Content.Find.Execute("688÷9=76, 4", $true, $false, $false, $false, $false, $true, 1, $false, "830÷5=166, 0", 2) | Out-Null
$d.Content.Find.Execute("899÷8=112, 3", $true, $false, $false, $false, $false, $true, 1, $false, "736÷3=245, 1", 2) | Out-Null
$d.Content.Find.Execute("969÷8=121, 1", $true, $false, $false, $false, $false, $true, 1, $false, "153÷7=21, 6", 2) | Out-Null
$d.Content.Find.Execute("424÷5=84, 4", $true, $false, $false, $false, $false, $true, 1, $false, "333÷8=41, 5", 2) | Out-Null
$d.Content.Find.Execute("467÷5=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "444÷6=74, 0", 2) | Out-Null
$d.Content.Find.Execute("381÷6=63, 3", $true, $false, $false, $false, $false, $true, 1, $false, "626÷2=313, 0", 2) | Out-Null
$d.Content.Find.Execute("837÷9=93, 0", $true, $false, $false, $false, $false, $true, 1, $false, "999÷8=124, 7", 2) | Out-Null
$d.Content.Find.Execute("343÷4=85, 3", $true, $false, $false, $false, $false, $true, 1, $false, "425÷6=70, 5", 2) | Out-Null
$d.Content.Find.Execute("559÷2=279, 1", $true, $false, $false, $false, $false, $true, 1, $false, "807÷7=115, 2", 2) | Out-Null
$d.Content.Find.Execute("759÷9=84, 3", $true, $false, $false, $false, $false, $true, 1, $false, "289÷9=32, 1", 2) | Out-Null
